$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 43
$ws.Cells.Item($row, 1).Value = "2025-04-29 04:34:04"
$ws.Cells.Item($row, 2).Value = 134
